# Word COM-interop script implementing:
#   "Added the option of specifying other course term plus added download
#    dropdown to syllabus views"
#
# Structural change inside word/document.xml:
#   - "${term}, ${season}" (one run) becomes two runs: "${term}" and
#     " ${season}", with the document's unique "_GoBack" bookmark moved to
#     sit right after "${term}" (between the two new runs). Word keeps only
#     a single "_GoBack" bookmark, so re-adding it at the new spot
#     automatically removes it from its old location (just before the
#     "${disability}" merge field near the end of the University Policies
#     section) and every "_Toc..." bookmark shifts up by one id, exactly as
#     in the target diff.

$d = $word.ActiveDocument

# 1) Turn "${term}, ${season}" into "${term} ${season}" (drop the comma).
$rng = $d.Content
$found = $rng.Find.Execute("`${term}, `${season}", $true, $false, $false, $false, $false, $true, 1, $false, "`${term} `${season}", 2)

# 2) Re-find "${term}" so we can collapse a range right after it, then drop
#    the (unique) "_GoBack" bookmark there -- Word automatically pulls
#    "_GoBack" off its previous location since only one can exist at a time.
$termRng = $d.Content
$termRng.Find.Execute("`${term}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$termRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $termRng)
